$wb = $excel.ActiveWorkbook

# --- Typography sheet: set Wildcard Characters (G4) and Wildcard Ranges (I4) to "0-9" ---
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTypo.Range("G4").Value = "0-9"
$wsTypo.Range("I4").Value = "0-9"

# --- Translation sheet: add two new text rows (4 and 5), and a trailing empty row (6) ---
$wsTrans = $wb.Worksheets.Item("Translation")

# Row 4: SingleUseId2 / Default / Left / LTR / "Value: <value>"
$wsTrans.Range("B4").Value = "SingleUseId2"
$wsTrans.Range("C4").Value = "Default"
$wsTrans.Range("D4").Value = "Left"
$wsTrans.Range("E4").Value = "LTR"
$wsTrans.Range("F4").Value = "Value: <value>"

# Row 5: SingleUseId3 / Default / Left / LTR / "0"
$wsTrans.Range("B5").Value = "SingleUseId3"
$wsTrans.Range("C5").Value = "Default"
$wsTrans.Range("D5").Value = "Left"
$wsTrans.Range("E5").Value = "LTR"
# Force "0" to be stored as text (shared string) instead of a number, then drop the
# text-number-format flag again so no extra style is left behind on the cell.
$wsTrans.Range("F5").NumberFormat = "@"
$wsTrans.Range("F5").Value = "0"
$wsTrans.Range("F5").Style = "Normal"

# Row 6: trailing empty row marker (matches the source row that was left blank)
$wsTrans.Rows(6).Hidden = $true
$wsTrans.Rows(6).Hidden = $false
